$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert a new row at position 14 (1-indexed), shifting existing rows 14-17 down to 15-18
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the new student's data
$ws.Cells.Item(14, 1).Value = 20330051920273
$ws.Cells.Item(14, 2).Value = "LUNA"
$ws.Cells.Item(14, 3).Value = "FLORES"
$ws.Cells.Item(14, 4).Value = "MIRANDA"
$ws.Cells.Item(14, 5).Value = "GEOMETRÍA Y TRIGONOMETRÍA"
$ws.Cells.Item(14, 6).Value = "2APV"
$ws.Cells.Item(14, 7).Value = 1
